$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.384.94"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.638.50"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.68%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.52"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.526"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.52%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.82"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.40%  "

$ws.Range("E9").Value = "  -2.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0609"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.98%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0892"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.870.31"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.69%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.635.39"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.02"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.560"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.26"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.354.27"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.49"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0719"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.53"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.74%  "

$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.30"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.60"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.22%  "

$ws.Range("E24").Value = "  -0.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.77"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.95"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.37%  "

$ws.Range("E27").Value = "  +0.86%  "

$ws.Range("E28").Value = "  +0.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.50"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.81%  "

$ws.Range("E30").Value = "  -4.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0482"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.90%  "

$ws.Range("E32").Value = "  -2.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.12"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.408.97"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.24%  "

$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("E36").Value = "  -0.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.562"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.878"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.51%  "

$ws.Range("E39").Value = "  -3.55%  "

$ws.Range("E40").Value = "  +0.90%  "

$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.47"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.74%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.49"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.82%  "

$ws.Range("E44").Value = "  +0.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.793"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.18"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -7.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.779.99"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.64"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.68"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.73%  "

$ws.Range("E50").Value = "  -1.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0984"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.51%  "

